# Add a new "UrunKodu" column to the "products" sheet.
#
# The authored edit inserted a new 4th data column ("UrunKodu", filled with
# "-") just before the existing "UrunAdi" column. Concretely:
#   - Header row: old D1 ("UrunAdi") text is moved to E1, D1 becomes "UrunKodu".
#   - Data rows 2-27: D keeps its original ("UrunAdi") value, and a new E
#     cell is added containing "-".
# Also: the "products" tab becomes the active/selected tab (it was
# "sections" before), with the selection left at E28 on "products" and
# H9 on "sections".

$wb = $excel.ActiveWorkbook

$wsProducts = $wb.Worksheets.Item("products")
$wsSections = $wb.Worksheets.Item("sections")

$headerD = $wsProducts.Cells.Item(1, 4)
$headerE = $wsProducts.Cells.Item(1, 5)

# Move the old header text ("UrunAdi") from D1 to the new E1 cell, then
# give D1 its new label.
$oldHeaderText = $headerD.Value2
$headerE.Value = $oldHeaderText
$headerD.Value = "UrunKodu"

# Copy the header's formatting (bold font, borders, centered alignment)
# onto the newly-populated E1 cell so it matches the rest of the header row.
$headerD.Copy()
$headerE.PasteSpecial(-4122)

# Fill the new column with "-" placeholders for every data row.
for ($r = 2; $r -le 27; $r++) {
    $wsProducts.Cells.Item($r, 5).Value = "-"
}

# "products" becomes the active sheet/tab, with E28 selected.
$wsProducts.Activate()
$null = $wsProducts.Range("E28").Select()

# "sections" keeps its own prior selection (H9); just make sure it is no
# longer the active tab (handled above by activating "products").
$null = $wsSections.Range("H9").Select()
$wsProducts.Activate()
